# Update countries & provincias Spain
# - Reorder "Angola" above "Sierra Leona"
# - Reorder "Eritrea" above "Islas Turcas y Caicos"
# - Update the "Datos actualizados" timestamp
# - Refresh several countries' COVID-19 statistics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 18 de Agosto de 2020 a las 22:03"

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 5638651
$ws.Cells.Item(4, 3).Value = 26624
$ws.Cells.Item(4, 4).Value = 2986346
$ws.Cells.Item(4, 5).Value = 2477707
$ws.Cells.Item(4, 7).Value = 882
$ws.Cells.Item(4, 8).Value = 174598

# Costa Rica (row 68)
$ws.Cells.Item(68, 2).Value = 29643
$ws.Cells.Item(68, 3).Value = 559
$ws.Cells.Item(68, 4).Value = 9462
$ws.Cells.Item(68, 5).Value = 19867
$ws.Cells.Item(68, 7).Value = 10
$ws.Cells.Item(68, 8).Value = 314

# Guayana Francesa (row 92)
$ws.Cells.Item(92, 2).Value = 8657
$ws.Cells.Item(92, 3).Value = 35
$ws.Cells.Item(92, 4).Value = 8054
$ws.Cells.Item(92, 5).Value = 550

# Albania (row 99)
$ws.Cells.Item(99, 2).Value = 7654
$ws.Cells.Item(99, 3).Value = 155
$ws.Cells.Item(99, 4).Value = 3871
$ws.Cells.Item(99, 5).Value = 3551
$ws.Cells.Item(99, 7).Value = 2
$ws.Cells.Item(99, 8).Value = 232

# Malaui (row 107)
$ws.Cells.Item(107, 2).Value = 5193
$ws.Cells.Item(107, 3).Value = 68
$ws.Cells.Item(107, 4).Value = 2716
$ws.Cells.Item(107, 5).Value = 2314
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 163

# Mali (row 126)
$ws.Cells.Item(126, 2).Value = 2666
$ws.Cells.Item(126, 3).Value = 26
$ws.Cells.Item(126, 4).Value = 1990
$ws.Cells.Item(126, 5).Value = 551

# Swap Sierra Leona / Angola (rows 137-138): Angola now comes first with
# refreshed stats, Sierra Leona follows with its previous (unchanged) stats.
$ws.Cells.Item(137, 1).Value = "Angola"
$ws.Cells.Item(137, 2).Value = 1966
$ws.Cells.Item(137, 3).Value = 31
$ws.Cells.Item(137, 4).Value = 667
$ws.Cells.Item(137, 5).Value = 1209
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 2
$ws.Cells.Item(137, 8).Value = 90

$ws.Cells.Item(138, 1).Value = "Sierra Leona"
$ws.Cells.Item(138, 2).Value = 1956
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 1506
$ws.Cells.Item(138, 5).Value = 381
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 69

# Move Eritrea above Islas Turcas y Caicos / Mongolia (rows 180-182):
# Eritrea now comes first with refreshed stats, Islas Turcas y Caicos and
# Mongolia follow, each keeping their previous (unchanged) stats.
$ws.Cells.Item(180, 1).Value = "Eritrea"
$ws.Cells.Item(180, 2).Value = 304
$ws.Cells.Item(180, 3).Value = 19
$ws.Cells.Item(180, 4).Value = 261
$ws.Cells.Item(180, 5).Value = 43
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(181, 2).Value = 298
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 55
$ws.Cells.Item(181, 5).Value = 241
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 2

$ws.Cells.Item(182, 1).Value = "Mongolia"
$ws.Cells.Item(182, 2).Value = 298
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 278
$ws.Cells.Item(182, 5).Value = 20
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0
